# Add a new "Jumlah Sampel" column (H) to the sample allocation sheet,
# mirroring the text style used by the existing header row, and populate
# the sample-count values for the two data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - same text-formatted style as the other header cells (e.g. G1)
$ws.Range("H1").Value = "Jumlah Sampel"
$ws.Range("H1").NumberFormat = "@"

# Data values (plain numbers, default style)
$ws.Range("H2").Value = 10
$ws.Range("H3").Value = 10

# Move/restore the active selection to match the saved view state
$null = $ws.Range("I10").Select()
